$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "72.040.90"
$cell.Style = $origStyle
$cell = $ws.Range("E2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.26%  "
$cell.Style = $origStyle
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.049.21"
$cell.Style = $origStyle
$cell = $ws.Range("E3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.04%  "
$cell.Style = $origStyle
$cell = $ws.Range("E4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.05%  "
$cell.Style = $origStyle
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "523.18"
$cell.Style = $origStyle
$cell = $ws.Range("E5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.43%  "
$cell.Style = $origStyle
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "148.64"
$cell.Style = $origStyle
$cell = $ws.Range("E6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.58%  "
$cell.Style = $origStyle
$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.624"
$cell.Style = $origStyle
$cell = $ws.Range("E7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.29%  "
$cell.Style = $origStyle
$cell = $ws.Range("E8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.10%  "
$cell.Style = $origStyle
$cell = $ws.Range("E9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.51%  "
$cell.Style = $origStyle
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.175"
$cell.Style = $origStyle
$cell = $ws.Range("E10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.97%  "
$cell.Style = $origStyle
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0000339"
$cell.Style = $origStyle
$cell = $ws.Range("E11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.10%  "
$cell.Style = $origStyle
$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "47.14"
$cell.Style = $origStyle
$cell = $ws.Range("E12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +10.83%  "
$cell.Style = $origStyle
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "10.71"
$cell.Style = $origStyle
$cell = $ws.Range("E13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.50%  "
$cell.Style = $origStyle
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.699.88"
$cell.Style = $origStyle
$cell = $ws.Range("E14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.17%  "
$cell.Style = $origStyle
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.054.90"
$cell.Style = $origStyle
$cell = $ws.Range("E15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.82%  "
$cell.Style = $origStyle
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "21.42"
$cell.Style = $origStyle
$cell = $ws.Range("E16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +7.74%  "
$cell.Style = $origStyle
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "14.32"
$cell.Style = $origStyle
$cell = $ws.Range("E17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.42%  "
$cell.Style = $origStyle
$cell = $ws.Range("E18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.69%  "
$cell.Style = $origStyle
$cell = $ws.Range("E19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.73%  "
$cell.Style = $origStyle
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "71.945.50"
$cell.Style = $origStyle
$cell = $ws.Range("E20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.41%  "
$cell.Style = $origStyle
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "441.95"
$cell.Style = $origStyle
$cell = $ws.Range("E21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.22%  "
$cell.Style = $origStyle
$cell = $ws.Range("B22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "Litecoin"
$cell.Style = $origStyle
$cell = $ws.Range("C22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$cell.Style = $origStyle
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "96.10"
$cell.Style = $origStyle
$cell = $ws.Range("E22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +8.98%  "
$cell.Style = $origStyle
$cell = $ws.Range("B23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "ImmutableX"
$cell.Style = $origStyle
$cell = $ws.Range("C23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell.Style = $origStyle
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.53"
$cell.Style = $origStyle
$cell = $ws.Range("E23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.79%  "
$cell.Style = $origStyle
$cell = $ws.Range("B24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "InternetComputer(DFINITY)"
$cell.Style = $origStyle
$cell = $ws.Range("C24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell.Style = $origStyle
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "14.46"
$cell.Style = $origStyle
$cell = $ws.Range("E24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.03%  "
$cell.Style = $origStyle
$cell = $ws.Range("B25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "RenderToken"
$cell.Style = $origStyle
$cell = $ws.Range("C25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell.Style = $origStyle
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.18"
$cell.Style = $origStyle
$cell = $ws.Range("E25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.85%  "
$cell.Style = $origStyle
$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.08"
$cell.Style = $origStyle
$cell = $ws.Range("E26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.16%  "
$cell.Style = $origStyle
$cell = $ws.Range("E27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +5.29%  "
$cell.Style = $origStyle
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "37.23"
$cell.Style = $origStyle
$cell = $ws.Range("E28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.58%  "
$cell.Style = $origStyle
$cell = $ws.Range("E29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.01%  "
$cell.Style = $origStyle
$cell = $ws.Range("B30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "Toncoin"
$cell.Style = $origStyle
$cell = $ws.Range("C30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$cell.Style = $origStyle
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.07"
$cell.Style = $origStyle
$cell = $ws.Range("E30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +7.92%  "
$cell.Style = $origStyle
$cell = $ws.Range("B31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "Bittensor"
$cell.Style = $origStyle
$cell = $ws.Range("C31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell.Style = $origStyle
$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "707.39"
$cell.Style = $origStyle
$cell = $ws.Range("E31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.40%  "
$cell.Style = $origStyle
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "13.51"
$cell.Style = $origStyle
$cell = $ws.Range("E32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.75%  "
$cell.Style = $origStyle
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.129"
$cell.Style = $origStyle
$cell = $ws.Range("E33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.01%  "
$cell.Style = $origStyle
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.89"
$cell.Style = $origStyle
$cell = $ws.Range("E34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +14.99%  "
$cell.Style = $origStyle
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "67.59"
$cell.Style = $origStyle
$cell = $ws.Range("E35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -4.10%  "
$cell.Style = $origStyle
$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0907"
$cell.Style = $origStyle
$cell = $ws.Range("E36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +6.90%  "
$cell.Style = $origStyle
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.444"
$cell.Style = $origStyle
$cell = $ws.Range("E37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -1.29%  "
$cell.Style = $origStyle
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "40.94"
$cell.Style = $origStyle
$cell = $ws.Range("E38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.37%  "
$cell.Style = $origStyle
$cell = $ws.Range("E39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +21.55%  "
$cell.Style = $origStyle
$cell = $ws.Range("E40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.21%  "
$cell.Style = $origStyle
$cell = $ws.Range("E41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.06%  "
$cell.Style = $origStyle
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = $origStyle
$cell = $ws.Range("E42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.19%  "
$cell.Style = $origStyle
$cell = $ws.Range("E43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.00%  "
$cell.Style = $origStyle
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.84"
$cell.Style = $origStyle
$cell = $ws.Range("E44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.94%  "
$cell.Style = $origStyle
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.12"
$cell.Style = $origStyle
$cell = $ws.Range("E45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.81%  "
$cell.Style = $origStyle
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.53"
$cell.Style = $origStyle
$cell = $ws.Range("E46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +5.05%  "
$cell.Style = $origStyle
$cell = $ws.Range("E47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.01%  "
$cell.Style = $origStyle
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.19"
$cell.Style = $origStyle
$cell = $ws.Range("E48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.05%  "
$cell.Style = $origStyle
$cell = $ws.Range("B49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "THORChain"
$cell.Style = $origStyle
$cell = $ws.Range("C49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$cell.Style = $origStyle
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.19"
$cell.Style = $origStyle
$cell = $ws.Range("E49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +7.55%  "
$cell.Style = $origStyle
$cell = $ws.Range("B50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "FLOKI"
$cell.Style = $origStyle
$cell = $ws.Range("C50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$cell.Style = $origStyle
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.000278"
$cell.Style = $origStyle
$cell = $ws.Range("E50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +22.14%  "
$cell.Style = $origStyle
$cell = $ws.Range("E51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.48%  "
$cell.Style = $origStyle
